# Fix resource-after-post issue: correct the Billing sheet's "id" value
# (was stale at 155, should be 11 to match the newly created resource),
# and leave the UI focused on the Billing sheet/cell instead of ProductDesc.

$wb = $excel.ActiveWorkbook

# --- Data fix -------------------------------------------------------------
$wsBilling = $wb.Worksheets.Item("Billing")
$wsBilling.Range("A2").Value = 11

# --- Selection / active-tab bookkeeping ------------------------------------
# Move ProductDesc's lingering selection back to A2 before switching away
# from it, then activate Billing and select its A2 cell so it becomes the
# workbook's active sheet/view.
$wsProduct = $wb.Worksheets.Item("ProductDesc")
$wsProduct.Activate()
$wsProduct.Range("A2").Select()

$wsBilling.Activate()
$wsBilling.Range("A2").Select()
